# Generate Report for Handoff
# Re-generates the handoff XLIFF for the "Ready for handoff" source files that
# have not yet been handed back (rows 4,6,7,8,9,10 in each locale sheet),
# refreshing the handoff timestamp/priority on the locale sheets and the
# "Latest HO Xliff Generate Date" on the Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(4, 6, 7, 8, 9, 10)

foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-10-13 14:19:32"

    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-10-13 14:19:20"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-10-13 14:19:32"
}
